# Update database and shift the reporting columns (read_price algorithm change):
# each period's figures move one column to the left (D<-E, E<-F, F<-G, G<-H)
# and the newest period (H) is populated with fresh data. The published-date
# header row (row 9) follows the same shift, with a brand new filing date in H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$cols = @("D","E","F","G","H")

# Row 9: "تاریخ انتشار" (publish date) header strip.
$row9 = @("1399-05-12 (5)", "1400-04-15 (9)", "1401-04-05 (10)", "1402-02-27 (7)", "1402-02-27")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "9").Value2 = $row9[$i]
}

# Data rows 11-27, each holding 5 period values (old col D..H -> new col D..H).
$data = @{
    11 = @(1091753, 1194322, 2247309, 5643364, 9761347)
    12 = @(-756513, -895999, -1759243, -4622082, -7181628)
    13 = @(335240, 298323, 488066, 1021282, 2579719)
    14 = @(-45866, -68836, -93011, -143320, -544332)
    15 = @(0, 0, 0, 0, 0)
    16 = @(0, 0, 1828, -147958, -221294)
    17 = @(289374, 229487, 396883, 730004, 1814093)
    18 = @(0, 0, 0, -158618, -694108)
    19 = @(26840, 42501, 42748, 165321, 86245)
    20 = @(316214, 271988, 439631, 736707, 1206230)
    21 = @(-37302, -16775, -23580, -57103, -66576)
    22 = @(278912, 255213, 416051, 679604, 1139654)
    23 = @(0, 0, 0, 0, 0)
    24 = @(278912, 255213, 416051, 679604, 1139654)
    25 = @(902, 365, 594, 971, 1628)
    26 = @(309166, 700000, 700000, 700000, 700000)
    27 = @(398, 365, 594, 971, 1628)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value2 = $vals[$i]
    }
}
